$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.146.26"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.024.88"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "226.56"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "55.13"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -5.85%  "
$ws.Range("D12").Value = "2.313.80"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "14.28"
$ws.Range("E13").Value = "  -4.75%  "
$ws.Range("D14").Value = "20.27"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "0.743"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "5.19"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "2.016.96"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "37.103.81"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").Value = "223.51"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -5.76%  "
$ws.Range("D26").Value = "166.49"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("D28").Value = "0.127"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "18.70"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "4.52"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "4.47"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("E35").Value = "  -4.87%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "5.58"
$ws.Range("E38").Value = "  +6.09%  "
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("D40").Value = "1.471.13"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("D42").Value = "95.88"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "16.38"
$ws.Range("E43").Value = "  -4.23%  "
$ws.Range("D44").Value = "0.0911"
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "1.14"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "2.75"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("D47").Value = "7.25"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "2.207.60"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "3.56"
$ws.Range("E51").Value = "  -10.33%  "
